$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI LR-pair values per Dr Hou's advice (ligand/receptor-expressing
# cell counts changed from 1 to 3, with corresponding recalculated stats).
# NOTE: PowerShell hashtables are case-insensitive, so the row-number key
# cannot be named "r" (it would collide with column "R"); use "RowNum".

$rows = @(
    @{ RowNum=2; E=3; G=2.147752666666667;  H=6.443258;            I=0.8708065646157738; J=0.8708065646157739;
       K=3; M=2.442429333333334;  N=7.327288;  O=0.08913295894744963; P=0.08913295894744963;
       Q=5.245734113811556; R=47.211607024304;  S=0.07761756577506741; T=0.07761756577506741 },
    @{ RowNum=3; E=3; G=2.147752666666667;  H=6.443258;            I=0.8708065646157738; J=0.8708065646157739;
       K=3; M=15.82990933333333; N=47.489728; O=0.5776898596383203;  P=0.5776898596383203;
       Q=33.99872998375822; R=305.988569853824; S=0.5030561220850143;  T=0.5030561220850143 },
    @{ RowNum=4; E=3; G=2.147752666666667;  H=6.443258;            I=0.8708065646157738; J=0.8708065646157739;
       K=3; M=9.129750999999999; N=27.389253; O=0.3331771814142301;  P=0.3331771814142301;
       Q=19.60844705625266; R=176.476023506274; S=0.2901328767556921;  T=0.2901328767556922 },
    @{ RowNum=5; E=3; G=0.318642;           H=0.9559260000000001; I=0.1291934353842261; J=0.1291934353842261;
       K=3; M=2.442429333333334;  N=7.327288;  O=0.08913295894744963; P=0.08913295894744963;
       Q=0.7782605676320001; R=7.004345108688001;  S=0.01151539317238221; T=0.01151539317238221 },
    @{ RowNum=6; E=3; G=0.318642;           H=0.9559260000000001; I=0.1291934353842261; J=0.1291934353842261;
       K=3; M=15.82990933333333; N=47.489728; O=0.5776898596383203;  P=0.5776898596383203;
       Q=5.044073969792001; R=45.396665728128;  S=0.074633737553306;   T=0.074633737553306 },
    @{ RowNum=7; E=3; G=0.318642;           H=0.9559260000000001; I=0.1291934353842261; J=0.1291934353842261;
       K=3; M=9.129750999999999; N=27.389253; O=0.3331771814142301;  P=0.3331771814142301;
       Q=2.909122118142;    R=26.182099063278;   S=0.04304430465853792; T=0.04304430465853792 }
)

foreach ($row in $rows) {
    $n = $row.RowNum
    $ws.Range("E$n").Value = $row.E
    $ws.Range("G$n").Value = $row.G
    $ws.Range("H$n").Value = $row.H
    $ws.Range("I$n").Value = $row.I
    $ws.Range("J$n").Value = $row.J
    $ws.Range("K$n").Value = $row.K
    $ws.Range("M$n").Value = $row.M
    $ws.Range("N$n").Value = $row.N
    $ws.Range("O$n").Value = $row.O
    $ws.Range("P$n").Value = $row.P
    $ws.Range("Q$n").Value = $row.Q
    $ws.Range("R$n").Value = $row.R
    $ws.Range("S$n").Value = $row.S
    $ws.Range("T$n").Value = $row.T
}
